$d = $word.ActiveDocument

# --- Equation 1: <ASE> -> <ASE>+<ASE>|<ASE>-<ASE>|<MDME>  becomes
#                 <ASE> -> <ASE>+<MDME>|<ASE>-<MDME>|<MDME>
# split across 5 m:r runs (same rPr as the original single run: Cambria Math, no size override)
$xml1 = '<m:oMathPara>' +
        '<m:oMath>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>&lt;ASE&gt; → &lt;ASE&gt;+&lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>MDM</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>E&gt;|&lt;ASE&gt;-&lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>MDM</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>E&gt;|&lt;MDME&gt;</m:t></m:r>' +
        '</m:oMath>' +
        '</m:oMathPara>'

$om1 = $d.OMaths.Item(1)
$om1.Range.InsertXML($xml1)

# --- Equation 2: <MDME> -> <MDME>*<MDME>|<MDME>/<MDME>|<MDME>%<MDME>|<PE>  becomes
#                 <MDME> -> <MDME>*<PE>|<MDME>/<PE>|<MDME>%<PE>|<PE>
# first run (<MDME> -> <) stays untouched; the second run is split into 7 m:r runs,
# all keeping the original run's rPr (Cambria Math, eastAsiaTheme minorEastAsia, sz 20)
$xml2 = '<m:oMathPara>' +
        '<m:oMath>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>&lt;MDME&gt; → &lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>MDME&gt;*&lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>P</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>E&gt;|&lt;MDME&gt;/&lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>P</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>E&gt;|&lt;MDME&gt;%&lt;</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>P</m:t></m:r>' +
          '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><m:t>E&gt;|&lt;PE&gt;</m:t></m:r>' +
        '</m:oMath>' +
        '</m:oMathPara>'

$om2 = $d.OMaths.Item(2)
$om2.Range.InsertXML($xml2)
